$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the range-prediction helper column G formula (F column values were
# originally captured in lakhs; they are corrected here to actual rupee
# amounts, and G keeps its previously-computed value as a static number
# instead of a live "=F*0.1" formula).
$data = @"
2,3900000.0,3.9000000000000004
3,6500000.0,6.5
4,7400000.0,7.4
5,8900000.0,8.9
6,7400000.0,7.4
7,9500000.0,9.5
8,7500000.0,7.5
9,10000000.0,10
10,6500000.0,6.5
11,8800000.0,8.8000000000000007
12,13500000.0,13.5
13,18800000.0,18.8
14,18800000.0,18.8
15,5800000.0,5.8000000000000007
16,9500000.0,9.5
17,4000000.0,4
18,5800000.0,5.8000000000000007
19,7800000.0,7.8000000000000007
20,3600000.0,3.6
21,4800000.0,4.8000000000000007
22,5700000.0,5.7
23,7900000.0,7.9
24,3700000.0,3.7
25,5700000.0,5.7
26,5900000.0,5.9
27,6100000.0,6.1000000000000005
28,6300000.0,6.3000000000000007
29,7900000.0,7.9
30,6100000.0,6.1000000000000005
31,6500000.0,6.5
32,7300000.0,7.3000000000000007
33,7600000.0,7.6000000000000005
34,8000000.0,8
35,9600000.0,9.6000000000000014
36,10000000.0,10
37,7000000.0,7
38,9000000.0,9
39,9700000.0,9.7000000000000011
40,12000000.0,12
41,4200000.0,4.2
42,5700000.0,5.7
43,TEXT,Null
44,1400000.0,1.4000000000000001
45,2200000.0,2.2000000000000002
46,2600000.0,2.6
47,8000000.0,8
48,8300000.0,8.3000000000000007
49,10100000.0,10.100000000000001
50,10300000.0,10.3
51,12800000.0,12.8
52,16300000.0,16.3
53,5400000.0,5.4
54,3900000.0,3.9000000000000004
55,5100000.0,5.1000000000000005
56,6000000.0,6
57,6600000.0,6.6000000000000005
58,4550000.0,4.55
59,7250000.0,7.25
60,7350000.0,7.3500000000000005
61,6200000.0,6.2
62,8400000.0,8.4
63,5414000.0,5.4140000000000006
64,5832000.0,5.8320000000000007
65,6938000.0,6.9379999999999997
66,7503000.0,7.5030000000000001
67,7902000.0,7.9020000000000001
68,8209999.999999999,8.2099999999999991
69,24900000.0,24.900000000000002
70,46000000.0,46
71,4575000.0,4.5750000000000002
72,4934000.0,4.9340000000000011
73,4974000.0,4.9740000000000002
74,5056000.0,5.0560000000000009
75,92300000.0,92.300000000000011
76,93000000.0,93
77,3600000.0,3.6
78,5500000.0,5.5
79,6550000.0,6.5500000000000007
80,6800000.0,6.8000000000000007
81,6500000.0,6.5
82,9200000.0,9.2000000000000011
83,9200000.0,9.2000000000000011
84,3350000.0,3.35
85,4730000.0,4.7299999999999995
86,5522000.0,5.5220000000000002
87,13600000.0,13.600000000000001
88,14300000.0,14.3
89,4759000.0,4.7590000000000003
90,4770000.0,4.7700000000000005
91,4780000.0,4.78
92,4791000.0,4.7909999999999995
93,4802000.0,4.8020000000000005
94,16100000.0,16.100000000000001
95,17000000.0,17
96,17800000.0,17.8
97,18700000.0,18.7
98,21100000.0,21.1
99,21300000.0,21.3
100,6500000.0,6.5
101,8400000.0,8.4
102,8500000.0,8.5
103,9600000.0,9.6000000000000014
104,12900000.0,12.9
105,13300000.0,13.3
106,13900000.0,13.9
107,14000000.0,14
108,4320000.0,4.32
109,4740000.0,4.74
110,3058000.0,3.0579999999999998
111,4349000.0,4.3490000000000002
112,4383000.0,4.383
113,4444000.0,4.444
114,4552000.0,4.5520000000000005
115,5400000.0,5.4
116,11100000.0,11.100000000000001
117,11800000.0,11.8
118,12100000.0,12.100000000000001
119,16200000.0,16.2
120,16700000.0,16.7
121,3600000.0,3.6
122,5200000.0,5.2
123,7900000.0,7.9
124,9700000.0,9.7000000000000011
125,6806000.0,6.8060000000000009
126,6841000.0,6.8410000000000002
127,7800000.0,7.8000000000000007
128,9887000.0,9.8870000000000005
129,14000000.0,14
130,17700000.0,17.7
131,6880000.0,6.88
132,6880000.0,6.88
133,13000000.0,13
134,4800000.0,4.8000000000000007
135,5100000.0,5.1000000000000005
136,5600000.0,5.6000000000000005
137,9500000.0,9.5
138,11000000.0,11
139,11800000.0,11.8
140,14600000.0,14.600000000000001
141,14700000.0,14.700000000000001
142,14900000.0,14.9
143,3868000.0,3.8680000000000003
144,5800000.0,5.8000000000000007
145,6072000.0,6.0720000000000001
146,8148999.999999999,8.1489999999999991
147,10800000.0,10.8
148,7450000.0,7.45
149,8059000.0,8.0590000000000011
150,8300000.0,8.3000000000000007
151,4725000.0,4.7250000000000005
152,4750000.0,4.75
153,4843000.0,4.843
154,5526000.0,5.5259999999999998
155,6600000.0,6.6000000000000005
156,6700000.0,6.7
157,4570000.0,4.57
158,4550000.0,4.55
159,5900000.0,5.9
160,8700000.0,8.7000000000000011
161,8900000.0,8.9
162,9200000.0,9.2000000000000011
163,4380000.0,4.38
164,5530000.0,5.53
165,5650000.0,5.65
166,6230000.0,6.23
167,7930000.0,7.93
168,8250000.0,8.25
169,15900000.0,15.9
170,21500000.0,21.5
171,4574000.0,4.5740000000000007
172,5895000.0,5.8950000000000005
173,5900000.0,5.9
174,10900000.0,10.9
175,3999000.0,3.9990000000000006
176,5399000.0,5.3990000000000009
177,5999000.0,5.9990000000000006
178,3300000.0,3.3000000000000003
179,4700000.0,4.7
180,7000000.0,7
181,7300000.0,7.3000000000000007
182,7800000.0,7.8000000000000007
183,8300000.0,8.3000000000000007
184,17500000.0,17.5
185,25000000.0,25
186,42000000.0,42
187,3949000.0,3.9490000000000003
188,3949000.0,3.9490000000000003
189,3949000.0,3.9490000000000003
190,6498999.999999999,6.4989999999999997
191,6498999.999999999,6.4989999999999997
192,6498999.999999999,6.4989999999999997
193,1319000.0,1.319
194,1729000.0,1.7290000000000001
195,8700000.0,8.7000000000000011
196,12900000.0,12.9
197,16000000.0,16
198,23000000.0,23
199,9500000.0,9.5
200,13500000.0,13.5
201,9068000.0,9.0680000000000014
"@

$rows = $data -split "`n" | Where-Object { $_.Trim() -ne "" }
foreach ($line in $rows) {
    $parts = $line.Trim() -split ","
    $r = [int]$parts[0]
    $fVal = $parts[1]
    $gVal = $parts[2]

    if ($fVal -eq "TEXT") {
        # F stays as-is (already the literal text, e.g. "Null");
        # G is set to the same text, replacing the old error formula.
        $ws.Cells.Item($r, 7).Value = $gVal
    } else {
        # Set G first so it becomes a plain cached number (no formula),
        # then update F to the corrected value.
        $ws.Cells.Item($r, 7).Value = [double]$gVal
        $ws.Cells.Item($r, 6).Value = [double]$fVal
    }
}

$null = $ws.Range("A2").Select()
